$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.722.73'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.27%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.456.02'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.29%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.46'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.69%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.48'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +8.15%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.457.08'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.36%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.05%  '

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.35%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.69'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.11%  '

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.34%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.393'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.67%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.049.79'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.40%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.97'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +10.09%  '

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.85%  '

# Row 16
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000175'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.72%  '

# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.470.43'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.52%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.828.46'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.14%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.26'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +9.07%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.38'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +4.10%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.59'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.85%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '388.02'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.05%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.566'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.75%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.44'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.52%  '

# Row 25
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.02%  '

# Row 26
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.78'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.33%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000123'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.80%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.606.48'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.50%  '

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.13%  '

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.93%  '

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.20%  '

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -8.01%  '

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.04%  '

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.28%  '

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.00%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '24.21'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.93%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.491.19'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.65%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.00'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.23%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.57'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.30%  '

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.31%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '166.62'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.41%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0788'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +4.44%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '27.19'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +7.26%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.809'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.58%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.52'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +4.07%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '42.52'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.98%  '

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.19%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.72'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.96%  '

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.88%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.575.38'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.79%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.94'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.31%  '
